# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets
# to reflect newly scraped totals (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6925
$ws1.Range("F3").Value = 402
$ws1.Range("F5").Value = 162
$ws1.Range("F6").Value = 6
$ws1.Range("F7").Value = 83
$ws1.Range("F8").Value = 593

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6925
$ws4.Range("F3").Value = 402
$ws4.Range("F6").Value = 162
$ws4.Range("F7").Value = 6
$ws4.Range("F9").Value = 83
$ws4.Range("F10").Value = 593
